$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: PREFIX value changed from "vocab" to "cl" ---
$ws.Range("B3").Value = "cl"

# --- Rows 20-40: update Identifier (A) and prefLabel (B); some rows gain a broader (F) value ---
$ws.Range("A20").Value = "cl:10000"
$ws.Range("B20").Value = "resource type "

$ws.Range("A21").Value = "cl:10001"
$ws.Range("B21").Value = "INSPIRE data theme"
$ws.Range("F21").Value = "cl:10000"

$ws.Range("A22").Value = "cl:10002"
$ws.Range("B22").Value = "unit"
$ws.Range("F22").Value = "cl:10000"

$ws.Range("A23").Value = "cl:10003"
$ws.Range("B23").Value = "spatial reference method"
$ws.Range("F23").Value = "cl:10000"

$ws.Range("A24").Value = "cl:10004"
$ws.Range("B24").Value = "coordinate reference system"

$ws.Range("A25").Value = "cl:10005"
$ws.Range("B25").Value = "administrative unit level"

$ws.Range("A26").Value = "cl:10006"
$ws.Range("B26").Value = "spatial scale"

$ws.Range("A27").Value = "cl:10007"
$ws.Range("B27").Value = "sampling time span"

$ws.Range("A28").Value = "cl:10008"
$ws.Range("B28").Value = "sampling time unit"

$ws.Range("A29").Value = "cl:10009"
$ws.Range("B29").Value = "method"

$ws.Range("A30").Value = "cl:10010"
$ws.Range("B30").Value = "data reliability"

$ws.Range("A31").Value = "cl:10011"
$ws.Range("B31").Value = "data completeness"

$ws.Range("A32").Value = "cl:10012"
$ws.Range("B32").Value = "processing step"

$ws.Range("A33").Value = "cl:10013"
$ws.Range("B33").Value = "resource provider"

$ws.Range("A34").Value = "cl:10014"
$ws.Range("B34").Value = "resource provider programme"

$ws.Range("A35").Value = "cl:10015"
$ws.Range("B35").Value = "geometry observation dataset"
$ws.Range("F35").Value = "cl:10000"

$ws.Range("A36").Value = "cl:10016"
$ws.Range("B36").Value = "statistical dataset"
$ws.Range("F36").Value = "cl:10000"

$ws.Range("A37").Value = "cl:10017"
$ws.Range("B37").Value = "measurement dataset"
$ws.Range("F37").Value = "cl:10000"

$ws.Range("A38").Value = "cl:10018"
$ws.Range("B38").Value = "category observation dataset"
$ws.Range("F38").Value = "cl:10000"

$ws.Range("A39").Value = "cl:10019"
$ws.Range("B39").Value = "count observation dataset"
$ws.Range("F39").Value = "cl:10000"

$ws.Range("A40").Value = "cl:10020"
$ws.Range("B40").Value = "truth observation dataset"
$ws.Range("F40").Value = "cl:10000"

# --- New row 41 ---
# Copy row 40's full A:Z range down to row 41 first so that every column
# (including the ones that stay blank) gets a real cell entry, matching the
# pattern used by all the other data rows on this sheet. Then overwrite the
# cells that actually carry new content.
$ws.Range("A40:Z40").Copy($ws.Range("A41:Z41"))
$ws.Range("A41").Value = "cl:10021"
$ws.Range("B41").Value = "temporal observation dataset"
$ws.Range("F41").Value = "cl:10000"
